# Apply updates to ValueSet-fr-medication-reconciliation-outcome workbook.
$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsInclude = $wb.Worksheets.Item("Include #0")

# URL: /fhir/fr/medication/ -> /ig/fhir/medication/
$wsMeta.Range("B2").Value = "https://hl7.fr/ig/fhir/medication/ValueSet/fr-medication-reconciliation-outcome"

# Name: FrMedicationReconciliationOutcome -> FRMedicationReconciliationOutcome
$wsMeta.Range("B4").Value = "FRMedicationReconciliationOutcome"

# Date updated
$wsMeta.Range("B8").Value = "2026-01-15T08:54:26+00:00"

# Jurisdiction value was empty, now set to FRANCE
$wsMeta.Range("B11").Value = "FRANCE"

# System URI on the Include sheet: /fhir/fr/medication/ -> /ig/fhir/medication/
$wsInclude.Range("B4").Value = "https://hl7.fr/ig/fhir/medication/CodeSystem/fr-medication-reconciliation-outcome"
